$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update labels: strip the trailing " is a problem" suffix
$ws.Range("A2").Value = "Income inequality in [Country]"
$ws.Range("A3").Value = "Climate change"
$ws.Range("A4").Value = "Global poverty"

# Overwrite figures with past ones with the right (more precise) size
$ws.Range("B2").Value = 0.550768403206271
$ws.Range("B3").Value = 0.591855289897272
$ws.Range("B4").Value = 0.502130679025185
